$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dados")

# Correct the mis-labeled "unnamed" header cells to "total"
$ws.Range("B2").Value = "total"
$ws.Range("F2").Value = "total"
